# Apply updated crypto price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.105.83"
$ws.Range("E2").Value = "  -0.30%  "

$ws.Range("D3").Value = "2.951.19"
$ws.Range("E3").Value = "  -1.29%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "569.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.46%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.55"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.37%  "

$ws.Range("E8").Value = "  +0.68%  "

$ws.Range("D9").Value = "2.945.36"
$ws.Range("E9").Value = "  -1.35%  "

$ws.Range("E10").Value = "  -4.73%  "

$ws.Range("E11").Value = "  -1.21%  "

$ws.Range("E12").Value = "  +2.06%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000244"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.87%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.06"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.62%  "

$ws.Range("E15").Value = "  -0.55%  "

$ws.Range("D16").Value = "65.152.27"
$ws.Range("E16").Value = "  -0.25%  "

$ws.Range("D17").Value = "3.440.20"
$ws.Range("E17").Value = "  -1.38%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.92"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.26%  "

$ws.Range("D19").Value = "2.952.73"
$ws.Range("E19").Value = "  -1.27%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "446.18"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.09%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.14"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.35%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.681"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.31%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.25"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.68%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.39"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.50%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.22"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.12%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.05"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.75%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.07"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.30%  "

$ws.Range("E28").Value = "  +0.07%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.70%  "

$ws.Range("E30").Value = "  -0.57%  "

$ws.Range("E31").Value = "  -0.58%  "

$ws.Range("E32").Value = "  -2.61%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.16"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.39%  "

$ws.Range("E34").Value = "  -0.88%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.15%  "

$ws.Range("E36").Value = "  -0.76%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.68"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.17%  "

$ws.Range("E38").Value = "  -0.31%  "

$ws.Range("E39").Value = "  -6.13%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "43.93"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.49%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.84"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.96%  "

$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.120"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.64%  "

$ws.Range("B43").Value = "TheGraph"
$ws.Range("C43").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.298"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.42%  "

$ws.Range("E44").Value = "  +0.10%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "387.08"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.13%  "

$ws.Range("E46").Value = "  +0.34%  "

$ws.Range("D47").Value = "2.718.28"
$ws.Range("E47").Value = "  -1.69%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "132.77"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.40%  "

$ws.Range("E49").Value = "  +0.11%  "

$ws.Range("E50").Value = "  +5.46%  "

$ws.Range("E51").Value = "  +0.71%  "
